$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows (170-179) appended to the Amazon sales trend sheet.
# Columns: A=SKU, B=Model, C=(Parent) ASIN, D=(Child) ASIN, P=units_ordered, T=ordered_product_sales
$rows = @(
    @("FBА" + "79266", "T-01",      "B0D2LGR3GG", "B0D2LGR3GG", 5, 4158.4799999999996),
    @("FBA79267",      "T-02",      "B0D2LD6BYJ", "B0D2LD6BYJ", 2, 1978.82),
    @("FBA79469",      "LE-04",     "B0DCK4DR1B", "B0DCK4DR1B", 1, 1232.2),
    @("FBA79332",      "HM-01",     "B0D63FKYZ5", "B0D63FKYZ5", 1, 862.86),
    @("FBA79346",      "ETC-04-WH", "B0D25LNDSY", "B0D25LNDSY", 0, 0),
    @("FBA79493",      "V-03",      "B0DQCWRG3H", "B0DQCWRG3H", 0, 0),
    @("FBA79271",      "V-01",      "B0D83Q7L8W", "B0D83Q7L8W", 0, 0),
    @("FBA79406",      "HSB-04",    "B0DCK3N2JJ", "B0DCK3N2JJ", 0, 0),
    @("FBA79464",      "CM-01-BL",  "B0DFCDKMWR", "B0DFCDKMWR", 0, 0),
    @("FBA79570",      "BR-01",     "B0DKJXRXKM", "B0DKJXRXKM", 0, 0)
)

$startRow = 170
$r = $startRow
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 16).Value = $row[4]
    $ws.Cells.Item($r, 20).Value = $row[5]
    $r = $r + 1
}
